# Apply the FormData sheet updates:
#  - remove the now-unused "timestamp" column (H)
#  - append new submission rows 19-27 in columns A-G

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column H ("timestamp" header + any data) entirely.
$ws.Columns.Item(8).Delete()

# New rows appended to the sheet (columns: lookingFor, planningToBuy, email,
# name, contactNumber, countryCode, readableTimestamp). Missing values are
# left blank, matching the sparse rows already present in the sheet.
$newRows = @(
    @("Life",   "1-2 Weeks", "ravikanttiwari488@gmail.com", "Ravikant Tiwari", "8744883594", "+40", "11/17/2024, 11:51:39 PM"),
    @("Life",   "1-2 Weeks", "ravikanttiwari488@gmail.com", "Ravikant Tiwari", "8744883594", "+40", "11/17/2024, 11:57:02 PM"),
    @("Travel", "1-2 Weeks", "ravikanttiwari488@gmail.com", "Ravikant Tiwari", "8744883594", "+40", "11/18/2024, 12:01:01 AM"),
    @("Life",   "1-2 Weeks", "ravikanttiwari488@gmail.com", "Ravikant Tiwari", "8744883594", "+40", "11/18/2024, 12:03:52 AM"),
    @("Life",   "1-2 Weeks", "ravikanttiwari488@gmail.com", "Ravikant Tiwari", "8744883594", "+40", "11/18/2024, 12:07:25 AM"),
    @("",       "1week",     "ravikanttiwari488@gmail.com", "Ravikant Tiwari", "8744883594", "+40", ""),
    @("Life",   "1-2 Weeks", "ravikanttiwari488@gmail.com", "Ravikant Tiwari", "8744883594", "+40", "11/18/2024, 12:27:26 AM"),
    @("",       "1week",     "rktindia2003@gmail.com",      "Ruchika kumari",  "9650511578", "+40", ""),
    @("",       "2week",     "ravikanttiwari488@gmail.com", "Ravikant Tiwari", "8744883594", "+40", "")
)

$startRow = 19
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $value = $rowData[$c]
        if ($value -ne "") {
            $cell = $ws.Cells.Item($r, $c + 1)
            if ($value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
                # Numeric-looking strings (phone numbers, the "+40" country
                # code) must stay literal text instead of being coerced into
                # numbers, so force text storage before assigning.
                $cell.NumberFormat = "@"
            }
            $cell.Value = $value
        }
    }
}
